$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "747÷6=124, 3" "520÷3=173, 1"
Replace-Text "301÷4=75, 1" "123÷4=30, 3"
Replace-Text "659÷9=73, 2" "637÷5=127, 2"
Replace-Text "713÷4=178, 1" "650÷4=162, 2"
Replace-Text "126÷5=25, 1" "655÷2=327, 1"
Replace-Text "987÷4=246, 3" "969÷5=193, 4"
Replace-Text "529÷8=66, 1" "826÷9=91, 7"
Replace-Text "154÷9=17, 1" "947÷9=105, 2"
Replace-Text "642÷8=80, 2" "539÷7=77, 0"
Replace-Text "795÷3=265, 0" "578÷3=192, 2"
Replace-Text "553÷4=138, 1" "260÷9=28, 8"
Replace-Text "968÷9=107, 5" "740÷7=105, 5"
Replace-Text "942÷2=471, 0" "234÷7=33, 3"
Replace-Text "159÷8=19, 7" "595÷9=66, 1"
Replace-Text "515÷2=257, 1" "372÷4=93, 0"
Replace-Text "145÷2=72, 1" "415÷4=103, 3"
Replace-Text "748÷9=83, 1" "363÷7=51, 6"
Replace-Text "417÷6=69, 3" "126÷9=14, 0"
Replace-Text "160÷2=80, 0" "243÷9=27, 0"
Replace-Text "845÷6=140, 5" "352÷7=50, 2"
Replace-Text "194÷8=24, 2" "629÷5=125, 4"
Replace-Text "279÷8=34, 7" "392÷4=98, 0"
Replace-Text "278÷7=39, 5" "102÷5=20, 2"
Replace-Text "214÷2=107, 0" "233÷6=38, 5"
Replace-Text "401÷4=100, 1" "548÷9=60, 8"

Write-Output "Done"
